$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying source records were re-sorted/re-exported: the data that was
# in rows 7-11 shifts up into rows 6-10, and the data that was in row 6 wraps
# around into row 11 (same six sighting records, new row order/IDs/coords).
# Apply the change cell-by-cell so that only the cells that actually differ
# between the old and new row order are touched, matching the diff exactly.

$ws.Range("A6").Value = 106607596
$ws.Range("Q6").Value = 405355.3536657504
$ws.Range("R6").Value = 7021553.610383645
$ws.Range("A7").Value = 106607595
$ws.Range("M7").Value = "färska spår"
$ws.Range("Q7").Value = 405340.130843634
$ws.Range("R7").Value = 7021495.582215455
$ws.Range("A8").Value = 106607577
$ws.Range("B8").Value = 89392
$ws.Range("E8").Value = 1202
$ws.Range("F8").Value = "Ullticka"
$ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("Q8").Value = 405210.1198294814
$ws.Range("R8").Value = 7021485.921191392
$ws.Range("AC8").ClearContents()
$ws.Range("A9").Value = 106607591
$ws.Range("B9").Value = 56395
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = "färska spår"
$ws.Range("N9").Value = ""
$ws.Range("Q9").Value = 405268.4341886034
$ws.Range("R9").Value = 7021309.663261802
$ws.Range("AC9").Value = "ringhack"
$ws.Range("A10").Value = 106607597
$ws.Range("M10").Value = ""
$ws.Range("Q10").Value = 405209.6556569744
$ws.Range("R10").Value = 7021470.190883989
$ws.Range("A11").Value = 106607590
$ws.Range("Q11").Value = 405235.6640926296
$ws.Range("R11").Value = 7021145.529926532
